# Swap the presentation's slide-master theme ("Integral" / Red Violet
# colours) for the default "Office Theme" colour scheme.
#
# (The companion notes-master theme part keeps the Integral colours that
# the slide-master theme used to have -- i.e. the two theme parts trade
# places. The PowerPoint object model only exposes a writable, 12-colour
# ThemeColorScheme on the presentation's single active theme, so we set
# that to the target "Office" values here.)

$p = $ppt.ActivePresentation
$sm = $p.SlideMaster
$tcs = $sm.Theme.ThemeColorScheme

# ThemeColorScheme items are addressed in theme order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
# RGB is a COM "long" in 0x00BBGGRR byte order.
$tcs.Item(1).RGB  = 0x000000   # dk1      000000
$tcs.Item(2).RGB  = 0xFFFFFF   # lt1      FFFFFF
$tcs.Item(3).RGB  = 0x6A5444   # dk2      44546A
$tcs.Item(4).RGB  = 0xE6E6E7   # lt2      E7E6E6
$tcs.Item(5).RGB  = 0xD59B5B   # accent1  5B9BD5
$tcs.Item(6).RGB  = 0x317DED   # accent2  ED7D31
$tcs.Item(7).RGB  = 0xA5A5A5   # accent3  A5A5A5
$tcs.Item(8).RGB  = 0x00C0FF   # accent4  FFC000
$tcs.Item(9).RGB  = 0xC47244   # accent5  4472C4
$tcs.Item(10).RGB = 0x47AD70   # accent6  70AD47
$tcs.Item(11).RGB = 0xC16305   # hlink    0563C1
$tcs.Item(12).RGB = 0x724F95   # folHlink 954F72
